$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Title" column (B) values: wrap the original title text in literal
# double quotes (the leading quote-prefix apostrophe is replaced by the
# quoted text starting with '=').
$ws.Range("B2").Value = "`"='Cur Ind (' & only(Year) & ')'`""
$ws.Range("B3").Value = "`"='Pr Sales (' & (only([Year]) - 1) & ')'`""
$ws.Range("B5").Value = "`"='Pr Profit (' & (only(Year) - 1) & ')'`""
$ws.Range("B7").Value = "`"='Pr Margin (' & (only(Year) - 1) & ')'`""
$ws.Range("B6").Value = "`"='Cur Margin (' & only(Year) & ')'`""
$ws.Range("B4").Value = "`"='Cur Profit (' & only(Year) & ')'`""

# These cells previously used a "quote prefix" style so the leading "'"
# would render literally. That's no longer needed now the text carries
# its own quoting, so reset the cell style back to Normal.
$ws.Range("B2:B7").Style = "Normal"

# Move the active selection to D3 (cosmetic cursor-position change).
[void]$ws.Range("D3").Select()
